# Update the cached "Date Placeholder" field text (footer date field) across
# every slide layout, the slide master, and the notes master from the old
# capture date (10/13/2021 / October 13, 2021) to the new one
# (3/18/2022 / March 18, 2022). Mirrors: "update repl links to use replit.com"
# commit's incidental re-save date bump.

function Update-DatePlaceholder($shp) {
    if ($shp.PlaceholderFormat.Type -eq 16) {
        $old = $shp.TextFrame.TextRange.Text
        if ($old -eq "10/13/2021") {
            $shp.TextFrame.TextRange.Text = "3/18/2022"
        } elseif ($old -eq "October 13, 2021") {
            $shp.TextFrame.TextRange.Text = "March 18, 2022"
        }
    }
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

# Slide master's own date placeholder.
foreach ($shp in $master.Shapes) {
    Update-DatePlaceholder $shp
}

# Every slide layout's date placeholder.
foreach ($layout in $master.CustomLayouts) {
    foreach ($shp in $layout.Shapes) {
        Update-DatePlaceholder $shp
    }
}

# Notes master's date placeholder.
$notesMaster = $p.NotesMaster
foreach ($shp in $notesMaster.Shapes) {
    Update-DatePlaceholder $shp
}
